# Update the MODS worksheet for the "Nosutu" item (liv_020018) and fix
# related quotation-mark styling on nearby entries in the "Documents" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Documents")

# Curly quote characters used in several titles below.
$lq = [char]0x201C   # "
$rq = [char]0x201D   # "

# --- Row 15: the volume number (93) was actually entered in the wrong
# column (Q = volume) and belongs in R (issue). Move it over, removing the
# now-empty Q15 cell entirely. ---
$ws.Range("R15").Value = 93
$ws.Range("Q15").Clear()

# --- Row 23 (new): add the new "Nosutu" item, a letter of gratitude from a
# Kafir woman (Nosutu Soga Jotelo), published in the Kaffrarian Watchman. ---
$ws.Range("A23").Value = "liv_020018"

# --- Row 21: "Jubilee of the Venerable Patriarch Brownlee" -> wrap title in
# curly quotes. ---
$ws.Range("B21").Value = $lq + "Jubilee of the Venerable Patriarch Brownlee" + $rq
$ws.Range("C21").Value = $lq + "Jubilee of the Venerable Patriarch Brownlee," + $rq + " 14 February 1867, 17 January 1867 "

# --- Row 20: "The Niger Expedtion" -> wrap title in curly quotes. ---
$ws.Range("B20").Value = $lq + "The Niger Expedtion" + $rq
$ws.Range("C20").Value = $lq + "The Niger Expedition," + $rq + " 25 September 1858, 24 June 1858 "

# --- Row 23 continued: title, alternative title, name, and genre. ---
$ws.Range("B23").Value = $lq + "Letter of Gratitude from a Kafir Woman" + $rq
$ws.Range("C23").Value = $lq + "Letter of Gratitude from a Kafir Woman," + $rq + " 1 May 1874"
$ws.Range("D23").Value = "Anonymous"
$ws.Range("E23").Value = "Jotelo, Nosutu Soga"
$ws.Range("G23").Value = " publications (documents)"
$ws.Rows.Item(23).RowHeight = 51

# --- Row 22: "Letter to A.M. Chirgwin" also needs a periodical citation
# added (journal, issue, pages, dates). ---
$ws.Range("P22").Value = "Wesleyan Juvenile Offering: A Miscellany of Missionary Info"
$ws.Range("R22").Value = 89
$ws.Range("S22").Value = "54, 56"

# --- Dates on rows 20 and 22 switch from comma- to semicolon-separated. ---
$ws.Range("T20").Value = "25 September 1858; 24 June 1858 "
$ws.Range("T22").Value = "1 May 1874; October 1873"

# Select the newly added cell, matching where the author was working.
$ws.Range("E23").Select()
